# Scheduled market-price refresh for the Leve-profit calculator sheets.
# Recomputed currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) per job sheet from the latest Universalis price snapshot.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 252.94737
$ws.Range("I39").Value = 81.38461
$ws.Range("J39").Value = 624.6667
$ws.Range("K39").Value = 244.15383
$ws.Range("L39").Value = 1874.0001
$ws.Range("M39").Value = 51.84617000000003
$ws.Range("N39").Value = -2466.0001
$ws.Range("H42").Value = 178.58824
$ws.Range("I42").Value = 142.4
$ws.Range("K42").Value = 427.2
$ws.Range("M42").Value = -197.2
$ws.Range("H76").Value = 4219.8
$ws.Range("J76").Value = 4498.4287
$ws.Range("L76").Value = 4498.4287
$ws.Range("N76").Value = -5128.4287
$ws.Range("H79").Value = 4219.8
$ws.Range("J79").Value = 4498.4287
$ws.Range("L79").Value = 4498.4287
$ws.Range("N79").Value = -6682.4287
$ws.Range("H92").Value = 408.3846
$ws.Range("I92").Value = 347.77777
$ws.Range("K92").Value = 347.77777
$ws.Range("M92").Value = 900.2222300000001
$ws.Range("H100").Value = 2090.6155
$ws.Range("I100").Value = 1668.4286
$ws.Range("K100").Value = 1668.4286
$ws.Range("M100").Value = -1127.4286
$ws.Range("H112").Value = 3115.25
$ws.Range("I112").Value = 2342.5715
$ws.Range("K112").Value = 7027.7145
$ws.Range("M112").Value = -5919.7145
$ws.Range("H131").Value = 3150.0667
$ws.Range("I131").Value = 1639.3077
$ws.Range("J131").Value = 12970
$ws.Range("K131").Value = 4917.9231
$ws.Range("L131").Value = 38910
$ws.Range("M131").Value = 122.0769
$ws.Range("N131").Value = -48990
$ws.Range("H135").Value = 1032.2941
$ws.Range("I135").Value = 516.8
$ws.Range("K135").Value = 4651.2
$ws.Range("M135").Value = -2116.2
$ws.Range("H141").Value = 1347.3823
$ws.Range("I141").Value = 1282.1515
$ws.Range("K141").Value = 3846.4545
$ws.Range("M141").Value = 1333.5455

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1856.3846
$ws.Range("I5").Value = 222.33333
$ws.Range("J5").Value = 3257
$ws.Range("K5").Value = 222.33333
$ws.Range("L5").Value = 3257
$ws.Range("M5").Value = -110.33333
$ws.Range("N5").Value = -3481
$ws.Range("H45").Value = 5163.8335
$ws.Range("I45").Value = 4495
$ws.Range("K45").Value = 4495
$ws.Range("M45").Value = -4118
$ws.Range("H61").Value = 8339037.5
$ws.Range("I61").Value = 6827.6665
$ws.Range("J61").Value = 33335666
$ws.Range("K61").Value = 6827.6665
$ws.Range("L61").Value = 33335666
$ws.Range("M61").Value = -6615.6665
$ws.Range("N61").Value = -33336090
$ws.Range("H63").Value = 1422.625
$ws.Range("I63").Value = 1384.1333
$ws.Range("K63").Value = 1384.1333
$ws.Range("M63").Value = -698.1333
$ws.Range("H66").Value = 1422.625
$ws.Range("I66").Value = 1384.1333
$ws.Range("K66").Value = 6920.666499999999
$ws.Range("M66").Value = -3488.666499999999
$ws.Range("H101").Value = 316730.22
$ws.Range("J101").Value = 316730.22
$ws.Range("L101").Value = 316730.22
$ws.Range("N101").Value = -323220.22
$ws.Range("H136").Value = 8339037.5
$ws.Range("I136").Value = 6827.6665
$ws.Range("J136").Value = 33335666
$ws.Range("K136").Value = 20482.9995
$ws.Range("L136").Value = 100006998
$ws.Range("M136").Value = -17932.9995
$ws.Range("N136").Value = -100012098

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1856.3846
$ws.Range("I4").Value = 222.33333
$ws.Range("J4").Value = 3257
$ws.Range("K4").Value = 222.33333
$ws.Range("L4").Value = 3257
$ws.Range("M4").Value = -107.33333
$ws.Range("N4").Value = -3487
$ws.Range("H22").Value = 1505
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1757.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1757.5
$ws.Range("M22").Value = -827
$ws.Range("N22").Value = -2103.5
$ws.Range("H35").Value = 35500
$ws.Range("J35").Value = 35500
$ws.Range("L35").Value = 35500
$ws.Range("N35").Value = -36120
$ws.Range("H86").Value = 3155.3125
$ws.Range("I86").Value = 3312.111
$ws.Range("K86").Value = 3312.111
$ws.Range("M86").Value = -2189.111
$ws.Range("H89").Value = 3155.3125
$ws.Range("I89").Value = 3312.111
$ws.Range("K89").Value = 16560.555
$ws.Range("M89").Value = -10944.555

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 790
$ws.Range("I2").Value = 487.5
$ws.Range("K2").Value = 487.5
$ws.Range("M2").Value = -374.5
$ws.Range("H7").Value = 1599.1177
$ws.Range("I7").Value = 31.666666
$ws.Range("K7").Value = 31.666666
$ws.Range("M7").Value = 81.33333400000001
$ws.Range("H17").Value = 431.25
$ws.Range("I17").Value = 431.25
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 431.25
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -257.25
$ws.Range("N17").ClearContents()
$ws.Range("H22").Value = 1420.28
$ws.Range("I22").Value = 177.72728
$ws.Range("J22").Value = 2396.5715
$ws.Range("K22").Value = 177.72728
$ws.Range("L22").Value = 2396.5715
$ws.Range("M22").Value = 172.27272
$ws.Range("N22").Value = -3096.5715
$ws.Range("H31").Value = 1357622.9
$ws.Range("I31").Value = 1427167.2
$ws.Range("K31").Value = 1427167.2
$ws.Range("M31").Value = -1426872.2
$ws.Range("H34").Value = 1357622.9
$ws.Range("I34").Value = 1427167.2
$ws.Range("K34").Value = 1427167.2
$ws.Range("M34").Value = -1426965.2
$ws.Range("H55").Value = 8333
$ws.Range("J55").Value = 9999
$ws.Range("L55").Value = 9999
$ws.Range("N55").Value = -10629

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 25121.75
$ws.Range("J49").Value = 25421
$ws.Range("L49").Value = 25421
$ws.Range("N49").Value = -25789
$ws.Range("H70").Value = 31946.525
$ws.Range("I70").Value = 30316.227
$ws.Range("J70").Value = 34188.188
$ws.Range("K70").Value = 30316.227
$ws.Range("L70").Value = 34188.188
$ws.Range("M70").Value = -30046.227
$ws.Range("N70").Value = -34728.188
$ws.Range("H73").Value = 31946.525
$ws.Range("I73").Value = 30316.227
$ws.Range("J73").Value = 34188.188
$ws.Range("K73").Value = 30316.227
$ws.Range("L73").Value = 34188.188
$ws.Range("M73").Value = -29380.227
$ws.Range("N73").Value = -36060.188
$ws.Range("H107").Value = 675.8333
$ws.Range("I107").Value = 619.5714
$ws.Range("K107").Value = 619.5714
$ws.Range("M107").Value = 1300.4286

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3336.5217
$ws.Range("I22").Value = 2075
$ws.Range("J22").Value = 3602.1052
$ws.Range("K22").Value = 2075
$ws.Range("L22").Value = 3602.1052
$ws.Range("M22").Value = -1780
$ws.Range("N22").Value = -4192.1052
$ws.Range("H27").Value = 3336.5217
$ws.Range("I27").Value = 2075
$ws.Range("J27").Value = 3602.1052
$ws.Range("K27").Value = 2075
$ws.Range("L27").Value = 3602.1052
$ws.Range("M27").Value = -1968
$ws.Range("N27").Value = -3816.1052
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H132").Value = 1590439.4
$ws.Range("I132").Value = 2383506.5
$ws.Range("K132").Value = 7150519.5
$ws.Range("M132").Value = -7147989.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 55095
$ws.Range("I70").Value = 55095
$ws.Range("K70").Value = 55095
$ws.Range("M70").Value = -54780
$ws.Range("H73").Value = 55095
$ws.Range("I73").Value = 55095
$ws.Range("K73").Value = 55095
$ws.Range("M73").Value = -54003
$ws.Range("H81").Value = 4267.15
$ws.Range("J81").Value = 4775
$ws.Range("L81").Value = 9550
$ws.Range("N81").Value = -11672
$ws.Range("H84").Value = 4267.15
$ws.Range("J84").Value = 4775
$ws.Range("L84").Value = 47750
$ws.Range("N84").Value = -58358
